$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.475.20"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.569.61"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.85"
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.501"
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.21"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0593"
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.791.45"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.572.94"
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("E15").Value = "  -2.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.81"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.468.09"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.59"
$ws.Range("E18").Value = "  -2.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0692"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.58"
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.01"
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.84"
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.72"
$ws.Range("E27").Value = "  -2.08%  "
$ws.Range("E28").Value = "  -0.60%  "
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0471"
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.379.61"
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("E34").Value = "  +1.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  +1.70%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.30"
$ws.Range("E36").Value = "  -0.88%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.951"
$ws.Range("E37").Value = "  -1.94%  "
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.546"
$ws.Range("E39").Value = "  +1.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.829"
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.979"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("E43").Value = "  +3.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.32"
$ws.Range("E44").Value = "  +1.27%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.27"
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.16"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.704.19"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.39"
$ws.Range("E48").Value = "  -3.27%  "
$ws.Range("E49").Value = "  -0.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0960"
$ws.Range("E50").Value = "  -1.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0497"
$ws.Range("E51").Value = "  -0.60%  "
